# Tracklist swap: Run the Jewels 3 -> Tyler, The Creator "Flower Boy"
# 1) Rename the two scoped defined names runthejewels3 -> tyler3
# 2) Replace the track data (title/composer/performer/time) on Sheet1 & Sheet3
# 3) Resize columns B/C/D on Sheet1 & Sheet3 to fit the new (longer) text
# Sheet2 is a pure formula "pretty-print" view of Sheet1, so it updates itself
# automatically on recalculation and needs no direct edits.

$wb = $excel.ActiveWorkbook

# --- 1. Defined names ------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws1.Names.Item(1).Name = "tyler3"
$ws3.Names.Item(1).Name = "tyler3"

# --- 2. New track data ------------------------------------------------------
# columns: B=Title  C=Composer  D=Performer  E=Time (fraction of a day)
$tracks = @(
    @("Foreword", "Michael Karoli / Jaki Liebezeit / Alex OConnor / Tyler Okonma / Irmin Schmidt / Holger Schuering / Damo Suzuki", "Tyler, The Creator feat. Rex Orange County", 0.13472222222222222),
    @("Where This Flower Blooms", "Frank Ocean / Tyler Okonma", "Tyler, The Creator feat. Frank ocean", 0.13472222222222222),
    @("Sometimes...", "Tyler Okonma", "Tyler, The Creator", 0.024999999999999998),
    @("See You Again", "Tyler Okonma", "Tyler, The Creator feat. Kali Uchis", 0.125),
    @("Who Dat Boy", "Rakim Mayers / Tyler Okonma", "Tyler, The Creator feat. A.$.A.P Rockey", 0.1423611111111111),
    @("Pothole", "Roy Ayers / Tyler Okonma", "Tyler, The Creator feat. Jaden Smith", 0.16388888888888889),
    @("Garden Shed", "Tyler Okonma / Estelle Swaray", "Tyler, The Creator feat. Estelle", 0.15486111111111112),
    @("Boredom", "Tyler Okonma", "Tyler, The Creator feat. Rex Orange County, Anna Of the North", 0.22222222222222221),
    @("I Ain't Got Time!", "Tyler Okonma", "Tyler, The Creator", 0.14305555555555557),
    @("911/Mr. Lonely", "Raymond Calhoun / Frank Ocean / Tyler Okonma", "Tyler, The Creator feat. Frank Ocean, Steve Lacy", 0.17708333333333334),
    @("Droppin' Seeds", "Dwayne Carter / Tyler Okonma", "Tyler, The Creator feat. Lil Wayne", 0.040972222222222222),
    @("November", "Tyler Okonma", "Tyler, The Creator", 0.15625),
    @("Glitter", "Tyler Okonma", "Tyler, The Creator", 0.15555555555555556),
    @("Enjoy Right Now Today", "Tyler Okonma", "Tyler, The Creator", 0.16319444444444445)
)

foreach ($ws in @($ws1, $ws3)) {
    $row = 2
    foreach ($t in $tracks) {
        $ws.Cells.Item($row, 2).Value = $t[0]
        $ws.Cells.Item($row, 3).Value = $t[1]
        $ws.Cells.Item($row, 4).Value = $t[2]
        $ws.Cells.Item($row, 5).Value = $t[3]
        $row++
    }

    # --- 3. Column widths --------------------------------------------------
    $ws.Columns.Item(2).ColumnWidth = 24.285714285714285
    $ws.Columns.Item(3).ColumnWidth = 80.14285714285714
    $ws.Columns.Item(4).ColumnWidth = 57.0
}
